$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leccion 4 (row 5) - English translation cell H5: join the two lines into
# a single paragraph (replace the embedded line break with a space).
$newText = "When generating the production packages, the environment settings were not correctly reviewed and certificates for that reason were sent with the test configuration and when trying to upload them to the cloud service it presented an error, which is why the customer complained in the attention given. Before generating a test or production package, the values of the keys must be verified within the Azure portal, in addition to this review the certificates. The web.config and cloud.ccproj files should be checked in the tags that say tests and / or production."

$ws.Range("H5").Value = $newText

# Update the active selection to match the saved workbook state.
$ws.Range("H6").Select()
